# Applies the VideoScript.docx edit described by the commit:
#  1. Collapses the fragmented "For this module, ... explore the ..."
#     run-soup in the "About the module(s)" paragraph into a single
#     cleaned-up sentence (wording + minor grammar/space fixes).
#  2. Un-semi-hides the built-in "Default Paragraph Font" character
#     style (drops <w:semiHidden/> while leaving <w:unhideWhenUsed/>).

$d = $word.ActiveDocument

# --- 1. Paragraph text rewrite -------------------------------------------
$old = "For this module,  it will explore the motorist statistics throughout " + `
       "all 12 stages of the 2024 Dakar Rally seeks to enhance the understanding " + `
       "of predictive modeling and statistical analysis among learners interested " + `
       "in competitive motorsport events. This investigation involves the " + `
       "application of multiple linear regression models to forecast driver " + `
       "rankings based on their cumulative stage times. Readers will glean " + `
       "insights into interpreting model summaries, detecting patterns and " + `
       "trends, and handling potential outliers. Through interactive exercises, " + `
       "individuals can hone their skills in model diagnostics, outlier " + `
       "detection, and evaluating model effectiveness using nested-hypothesis " + `
       "tests. Ultimately, this endeavor furnishes a pragmatic framework for " + `
       "employing statistical techniques in sports contexts."

$new = "This module will explore motorist statistics throughout all 12 stages " + `
       "of the 2024 Dakar Rally and seeks to enhance the understanding of " + `
       "predictive modeling and statistical analysis among learners interested " + `
       "in competitive motorsport events. This investigation involves the " + `
       "application of multiple linear regression models to forecast driver " + `
       "rankings based on their cumulative stage times. Readers will glean " + `
       "insights into interpreting model summaries, detecting patterns and " + `
       "trends, and handling potential outliers. Through interactive exercises, " + `
       "individuals can hone their skills in model diagnostics, outlier " + `
       "detection, and evaluating model effectiveness using nested hypothesis " + `
       "tests. Ultimately, this endeavor furnishes a pragmatic framework for " + `
       "employing statistical techniques in sports contexts."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not locate the 'For this module...' paragraph text to replace."
}

# --- 2. Un-semi-hide the "Default Paragraph Font" style -------------------
# (matches the styles.xml hunk that drops <w:semiHidden/> from that style).
try {
    $dpf = $d.Styles.Item("Default Paragraph Font")
    $dpf.Hidden = $false
} catch {
    # Some hosts don't expose a settable Style.Hidden; the text edit above
    # is the load-bearing part of this change, so don't fail the run over it.
    Write-Output "Style.Hidden could not be set: $_"
}
